# "Generate Report for Archive"
# - Localization status moves from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each
#   language sheet's own Status column.
# - The Status columns are narrower now that the new text is shorter,
#   matching the regenerated report's column sizing.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text wherever it appears ---

# Overview sheet: E2 = zh-cn status, F2 = de-de status
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }

# Language sheets: column C is "Status"
if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }

# --- Resize the Status columns to fit the new, shorter text ---
# (Target stored width ~13.41 chars; 12.5 is what the host snaps to that cell.)

$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = $newWidth   # F: de-de status

$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # C: Status
$dede.Columns.Item(3).ColumnWidth = $newWidth        # C: Status
